$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 215-216; everything from row 215 down shifts
# down by two rows (215->217 ... 324->326), matching the target diff.
$ws.Range("A215:A216").EntireRow.Insert()

# New row 215: Cilantro, Primera, $/caja 36 atados, Región Metropolitana
$ws.Cells.Item(215, 1).Value = 9
$ws.Cells.Item(215, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(215, 3).Value = "Metropolitana"
$ws.Cells.Item(215, 4).Value = 44455
$ws.Cells.Item(215, 5).Value = 13
$ws.Cells.Item(215, 6).Value = 100112040
$ws.Cells.Item(215, 7).Value = "Cilantro"
$ws.Cells.Item(215, 8).Value = "Sin especificar"
$ws.Cells.Item(215, 9).Value = "Primera"
$ws.Cells.Item(215, 10).Value = 52
$ws.Cells.Item(215, 11).Value = 4000
$ws.Cells.Item(215, 12).Value = 4500
$ws.Cells.Item(215, 13).Value = 4250
$ws.Cells.Item(215, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(215, 15).Value = "Región Metropolitana"
$ws.Cells.Item(215, 16).Value = 118
$ws.Cells.Item(215, 17).Value = 36
$ws.Cells.Item(215, 18).Value = "Hortaliza"

# New row 216: Cilantro, Primera, $/docena de atados, Región Metropolitana
$ws.Cells.Item(216, 1).Value = 9
$ws.Cells.Item(216, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(216, 3).Value = "Metropolitana"
$ws.Cells.Item(216, 4).Value = 44455
$ws.Cells.Item(216, 5).Value = 13
$ws.Cells.Item(216, 6).Value = 100112040
$ws.Cells.Item(216, 7).Value = "Cilantro"
$ws.Cells.Item(216, 8).Value = "Sin especificar"
$ws.Cells.Item(216, 9).Value = "Primera"
$ws.Cells.Item(216, 10).Value = 133
$ws.Cells.Item(216, 11).Value = 8000
$ws.Cells.Item(216, 12).Value = 10000
$ws.Cells.Item(216, 13).Value = 9008
$ws.Cells.Item(216, 14).Value = "`$/docena de atados"
$ws.Cells.Item(216, 15).Value = "Región Metropolitana"
$ws.Cells.Item(216, 16).Value = 3003
$ws.Cells.Item(216, 17).Value = 3
$ws.Cells.Item(216, 18).Value = "Hortaliza"
